$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cell from "Classroom " to "Name"
$ws.Range("A1").Value = "Name"

# Update the active selection to A2
$ws.Range("A2").Select()
